# Fix import errors and data processing issues
# - Replace the Dilemmas sheet's 'Dilemma' column (F) for rows 15-86 with
#   newly generated, more nuanced dilemma phrasings (mirrors a corrected
#   mistral/data-processing pass that produced richer dilemma text).
# - Restore UI selection state: Dilemmas becomes the active tab (selection
#   on E6), Dark_Triad_Results selection moves to B9 and is no longer the
#   active tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dilemmas")

$newDilemmas = @(
    'Should companies prioritise environmental responsibility even if it significantly raises consumer costs and impacts low-income consumers?'
    'Should economic stability and job security be risked to enforce strict environmental protections on industries vital to the local economy?'
    'Should personal freedom be restricted by mandating carbon footprint reductions for the collective good of future generations?'
    'Should the complete banning of deforestation be enforced even if it severely limits economic development and impacts communities dependent on forestry?'
    'Should endangered species be given legal rights even if it severely restricts human activities, potentially affecting local economies and communities?'
    'Should user data be freely traded if it improves services and personalisation but compromises personal privacy?'
    'Should employees sometimes accept unpaid overtime if refusing it jeopardises the company''s survival and their own employment?'
    'Should companies be globally accountable for wages, even if increased fairness dramatically raises costs for consumers and threatens the company''s survival?'
    'Should whistleblowers always be protected, even if their revelations significantly harm national security or public safety?'
    'Should businesses prioritise profit to sustain economic growth, even if this occasionally undermines social responsibilities and community welfare?'
    'Should euthanasia be permitted to respect personal choice, despite risks of misuse or moral concerns from society?'
    'Should organ donation be mandatory if it saves many lives, despite potentially infringing on personal autonomy after death?'
    'Should genetic editing of embryos be allowed to eliminate diseases, despite unknown risks and ethical concerns about ''designer babies''?'
    'Should patients have the right to refuse life-saving treatments even if their decision adversely impacts their families and healthcare resources?'
    'Should pharmaceutical companies freely set drug prices to fund research, even if many patients cannot afford life-saving medication?'
    'Should human cloning be permitted for scientific advancement, despite ethical concerns regarding identity and exploitation?'
    'Should insurers use genetic data for risk assessment, improving accuracy but potentially leading to discrimination and inequality?'
    'Should stem cell research receive federal funding to promote medical breakthroughs, despite ethical objections from segments of society?'
    'Should animal testing for medical research continue, given the ethical cost versus the potential for life-saving treatments?'
    'Should age restrictions on reproductive technologies balance societal stability against individual rights to family planning?'
    'Should the death penalty be abolished if it risks reducing deterrence of serious crimes, despite ethical objections to capital punishment?'
    'Should juveniles be tried as adults for severe crimes, prioritising justice for victims, despite potentially neglecting juvenile rehabilitation?'
    'Should prisons emphasise rehabilitation even if it means lighter sentences for serious crimes, risking public safety perceptions?'
    'Should mandatory minimum sentences exist to ensure consistent justice, even if they occasionally lead to disproportionate punishments?'
    'Should restorative justice be prioritised over traditional punishments, despite potential perceptions of leniency from victims?'
    'Should affirmative action policies be implemented, potentially disadvantaging individuals based on merit alone?'
    'Should universal healthcare be provided regardless of cost, potentially limiting resources in other vital areas like education or infrastructure?'
    'Should enforced wealth redistribution address inequality, despite potentially discouraging individual effort and innovation?'
    'Should equal pay across all jobs be mandated, potentially reducing motivation for skilled and demanding professions?'
    'Should political quotas ensure minority representation, possibly compromising meritocratic selection processes?'
    'Should gender equality laws mandate equal pay, even if it impacts companies'' flexibility in rewarding experience or performance?'
    'Should marriage equality be universally recognised despite cultural or religious objections?'
    'Should universal basic income be guaranteed, despite potential impacts on motivation, productivity, and economic sustainability?'
    'Should educational resources be uniformly distributed, potentially limiting resources to high-performing regions to support struggling ones?'
    'Should legal mandates for disability rights prioritise equity even if the financial burden impacts overall public funding?'
    'Should freedom of speech include limitations to prevent harm, potentially suppressing dissenting views or unpopular opinions?'
    'Should intervention to stop human rights abuses be required even if it means violating national sovereignty and risking wider conflict?'
    'Should privacy be absolute, potentially hindering security and criminal investigations?'
    'Should refugees be granted the same rights as citizens, potentially straining resources and creating tensions within host communities?'
    'Should internet access be a guaranteed right despite the significant costs involved, possibly diverting funds from essential services?'
    'Should vaccination refusal be permitted despite risks to public health and community safety?'
    'Should press freedom be unrestricted, even if misinformation could lead to public harm?'
    'Should freedom of expression protect hate speech, risking societal division and potential harm?'
    'Should euthanasia be a personal freedom despite potential ethical slippery slopes and societal discomfort?'
    'Should personal freedoms be restricted to maintain public safety, potentially risking authoritarian overreach?'
    'Should parents face legal consequences for children''s crimes, potentially punishing them unfairly for circumstances beyond their control?'
    'Should corporations be compelled to clean up environmental damage, even if costs threaten their economic viability and employment levels?'
    'Should governments be accountable for public health, even if it significantly expands government authority and costs?'
    'Should individuals solely manage retirement savings, despite risks of widespread poverty in old age?'
    'Should manufacturers be fully accountable for product safety, even if extensive regulation stifles innovation and increases consumer costs?'
    'Should politicians disclose financial interests despite potential invasion of privacy and increased vulnerability to political attacks?'
    'Should companies fully disclose supply chains despite potentially risking trade secrets and competitive advantage?'
    'Should social media platforms actively fact-check content, potentially limiting free expression or becoming gatekeepers of information?'
    'Should employees trust their employers with personal data, despite possible exploitation or breaches?'
    'Should scientific findings be trusted without independent verification, even if verification slows progress?'
    'Should journalists be compelled to reveal sources, potentially endangering future informants and press freedom?'
    'Should cheating always result in expulsion, despite potential harsh consequences for a student''s future?'
    'Should companies disclose all risks associated with products, potentially harming their competitive edge and consumer confidence?'
    'Should politicians resign due to personal scandals despite their effective public performance?'
    'Should researchers who falsify data face permanent bans, despite potential redemption and valuable expertise?'
    'Should employees prioritise company interests over personal ethics, potentially compromising individual integrity?'
    'Should friends report each other''s illegal activities, risking personal relationships for legal or moral duty?'
    'Should soldiers disobey orders perceived as unethical, despite risking discipline and military cohesion?'
    'Should consumers remain loyal to brands despite better alternatives, prioritising relationship over value?'
    'Should loyalty to one''s country override personal disagreement with policies, risking complicity in unethical practices?'
    'Should harmful cultural practices be respected to preserve cultural identity despite clear individual harm?'
    'Should parental authority be absolute, potentially overriding individual autonomy and moral judgement?'
    'Should all opinions be respected equally, even those clearly harmful or dangerous?'
    'Should individual privacy take precedence even if it obstructs significant public interest matters?'
    'Should elderly receive preferential treatment even if it disadvantages younger individuals?'
    'Should patients always control treatment decisions even against professional medical advice?'
    'Should children independently make major decisions despite risks of immature judgement?'
)

$startRow = 15
for ($i = 0; $i -lt $newDilemmas.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 6).Value = $newDilemmas[$i]
}

# Dark_Triad_Results: move the selection to B9 (no longer the active tab).
$wsDarkTriadResults = $wb.Worksheets.Item("Dark_Triad_Results")
$wsDarkTriadResults.Range("B9").Select()

# Dilemmas becomes the active tab, selection on E6.
$ws.Activate()
$ws.Range("E6").Select()

"Updated $($newDilemmas.Length) Dilemma cells (F15:F86)."
